# CIROH_Items_Comprehensive_List.xlsx update
# "Updated CIROH_Items_Comprehensive_List. Should now contain all orders from fall 2024"
#
# Adds new equipment rows (11-13, 15, 17-20) to the "2024-2025" sheet and
# stamps the existing order rows (6-9) with a "Fall 2024" date-received note
# in column I.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024-2025")
$ws.Activate()

$CURRENCY2 = '"$"#,##0.00_);[Red]("$"#,##0.00)'
$CURRENCY0 = '"$"#,##0_);[Red]("$"#,##0)'

# ---------------------------------------------------------------------
# Pre-seed the shared-string table in the same order the new strings were
# introduced, using scratch cells far away from the printed area. The
# cells get overwritten later with their real, final content.
# ---------------------------------------------------------------------
$seedRow = 500
$seedStrings = @(
    "Transducer Contact Microphone",
    "Amazon.com: DZS Elec 15PCS 35mm Piezo Disc Transducer Contact Microphone Trigger Sound Sensor with 4 Inches Wires for Acoustic Instrument : Musical Instruments",
    "CME0303S3C - Isolated PS",
    "Digikey",
    "Automatic gain controll Microphone",
    "Cable (50 feet)",
    "Microphone ",
    "Molex 51021-1000 ",
    "DigiKey",
    "Raspberry pi ",
    "Fall 2024",
    "Jordan",
    "RockBlock Tranceiver"
)
for ($i = 0; $i -lt $seedStrings.Count; $i++) {
    $ws.Cells.Item($seedRow, $i + 1).Value = $seedStrings[$i]
}

# ---------------------------------------------------------------------
# Existing rows 6-9: stamp column I with "Fall 2024"
# ---------------------------------------------------------------------
$ws.Range("I6").Value = "Fall 2024"
$ws.Range("I7").Value = "Fall 2024"
$ws.Range("I8").Value = "Fall 2024"
$ws.Range("I9").Value = "Fall 2024"

# ---------------------------------------------------------------------
# Row 11 - Transducer Contact Microphone
# ---------------------------------------------------------------------
$ws.Range("A11").Value = "Transducer Contact Microphone"
$ws.Range("A11").HorizontalAlignment = -4108
$ws.Range("A11").WrapText = $true

$ws.Range("B11").Value = 1
$ws.Range("B11").NumberFormat = $CURRENCY2
$ws.Range("B11").HorizontalAlignment = -4108

$ws.Range("C11").Value = 4.99
$ws.Range("C11").NumberFormat = $CURRENCY2
$ws.Range("C11").HorizontalAlignment = -4108

$ws.Range("D11").Value = 4.99
$ws.Range("D11").NumberFormat = $CURRENCY2

$ws.Range("E11").Value = "Amazon.com: DZS Elec 15PCS 35mm Piezo Disc Transducer Contact Microphone Trigger Sound Sensor with 4 Inches Wires for Acoustic Instrument : Musical Instruments"
$ws.Range("E11").WrapText = $true
$ws.Hyperlinks.Add($ws.Range("E11"), "https://www.amazon.com/DZS-Elec-Transducer-Microphone-Instrument/dp/B084KHH7B6/ref=sr_1_19?dib=eyJ2IjoiMSJ9.-PibDNlinRz79bIpdumGJKLNDqgwo8MJ5Yur3e9XNYKUOlqtMaPkT5RbNL1UTHdk9PlMjYh8k938pDwecAJ3MBVB5GgJp0aDSe5PK_QLgmFmKujHDdR2CkbwuTbrQngjvjEtKLfrqUjdOWQv28gJyphTne3z90EQrXzwTVcDgEgUVL52XbIaM1IbMUQgsKnKgYSKWdKSTIkeM-uaM6fUkM3cPWOA0B2ba7SIURtCjegOvJMqLDOWVxTC8Gz3ygrqHVCl8lisiaulh38ciug__mY2EYUkTgSYbwo3FhV8imc.8H3H4oIOGSwKLzykTkDj1mJYWcn7OKWwjFTXtMVt8R8&dib_tag=se&keywords=contact%2Bmicrophone&qid=1726071396&sr=8-19&th=1") | Out-Null

$ws.Range("I11").Value = "Fall 2024"

$ws.Rows.Item(11).RowHeight = 100.8

# ---------------------------------------------------------------------
# Row 12 - CME0303S3C - Isolated PS
# ---------------------------------------------------------------------
$ws.Range("A12").Value = "CME0303S3C - Isolated PS"
$ws.Range("A12").HorizontalAlignment = -4108
$ws.Range("A12").WrapText = $true

$ws.Range("B12").Value = 2
$ws.Range("B12").HorizontalAlignment = -4108

$ws.Range("C12").Value = 3.25
$ws.Range("C12").NumberFormat = $CURRENCY2
$ws.Range("C12").HorizontalAlignment = -4108

$ws.Range("D12").Value = 6.5
$ws.Range("D12").NumberFormat = $CURRENCY2

$ws.Range("E12").Value = "Digikey"
$ws.Hyperlinks.Add($ws.Range("E12"), "https://www.digikey.com/en/products/detail/cui-inc/CME0303S3C/7705660") | Out-Null
$ws.Range("E12").HorizontalAlignment = -4108

$ws.Range("I12").Value = "Fall 2024"

$ws.Rows.Item(12).RowHeight = 28.8

# ---------------------------------------------------------------------
# Row 13 - Automatic gain control Microphone
# ---------------------------------------------------------------------
$ws.Range("A13").Value = "Automatic gain controll Microphone"
$ws.Range("A13").HorizontalAlignment = -4108
$ws.Range("A13").WrapText = $true

$ws.Range("B13").Value = 1
$ws.Range("B13").HorizontalAlignment = -4108

$ws.Range("C13").Value = 7.95
$ws.Range("C13").NumberFormat = $CURRENCY2
$ws.Range("C13").HorizontalAlignment = -4108

$ws.Range("D13").Value = 7.95
$ws.Range("D13").NumberFormat = $CURRENCY2

$ws.Range("E13").Value = "Digikey"
$ws.Hyperlinks.Add($ws.Range("E13"), "https://www.digikey.com/en/products/detail/adafruit-industries-llc/1713/5273713") | Out-Null
$ws.Range("E13").HorizontalAlignment = -4108

$ws.Range("I13").Value = "Fall 2024"

$ws.Rows.Item(13).RowHeight = 28.8

# ---------------------------------------------------------------------
# Row 15 - RockBlock Tranceiver note / Jordan
# ---------------------------------------------------------------------
$ws.Range("A15").Value = "RockBlock Tranceiver"
$ws.Range("A15").HorizontalAlignment = -4108

$ws.Range("G15").Value = "Jordan"

$ws.Range("I15").Value = "Fall 2024"

# Row 16 - blank spacer row, column I keeps the copied-down (blank) style
$ws.Range("I16").Value = "Fall 2024"
$ws.Range("I16").Value = ""

# ---------------------------------------------------------------------
# Row 17 - Cable (50 feet)
# ---------------------------------------------------------------------
$ws.Range("A17").Value = "Cable (50 feet)"
$ws.Range("A17").HorizontalAlignment = -4108

$ws.Range("B17").Value = 1
$ws.Range("B17").HorizontalAlignment = -4108

$ws.Range("C17").Value = 21.99
$ws.Range("C17").NumberFormat = $CURRENCY2
$ws.Range("C17").HorizontalAlignment = -4108

$ws.Range("D17").Value = 21.99
$ws.Range("D17").NumberFormat = $CURRENCY2
$ws.Range("D17").HorizontalAlignment = -4108

$ws.Range("E17").Value = "Amazon"
$ws.Hyperlinks.Add($ws.Range("E17"), "https://www.amazon.com/s?k=cable+50+feet") | Out-Null
$ws.Range("E17").HorizontalAlignment = -4108

$ws.Range("I17").Value = "Fall 2024"

# ---------------------------------------------------------------------
# Row 18 - Microphone
# ---------------------------------------------------------------------
$ws.Range("A18").Value = "Microphone "
$ws.Range("A18").HorizontalAlignment = -4108

$ws.Range("B18").Value = 1
$ws.Range("B18").HorizontalAlignment = -4108

$ws.Range("C18").Value = 7.49
$ws.Range("C18").NumberFormat = $CURRENCY2
$ws.Range("C18").HorizontalAlignment = -4108

$ws.Range("D18").Value = 7.49
$ws.Range("D18").NumberFormat = $CURRENCY2
$ws.Range("D18").HorizontalAlignment = -4108

$ws.Range("E18").Value = "Amazon"
$ws.Hyperlinks.Add($ws.Range("E18"), "https://www.amazon.com/s?k=microphone") | Out-Null
$ws.Range("E18").HorizontalAlignment = -4108

$ws.Range("I18").Value = "Fall 2024"

# ---------------------------------------------------------------------
# Row 19 - Molex 51021-1000
# ---------------------------------------------------------------------
$ws.Range("A19").Value = "Molex 51021-1000 "
$ws.Range("A19").HorizontalAlignment = -4108

$ws.Range("B19").Value = 5
$ws.Range("B19").HorizontalAlignment = -4108

$ws.Range("C19").Value = 0.42
$ws.Range("C19").NumberFormat = $CURRENCY2
$ws.Range("C19").HorizontalAlignment = -4108

$ws.Range("D19").Value = 2.1
$ws.Range("D19").NumberFormat = $CURRENCY2
$ws.Range("D19").HorizontalAlignment = -4108

$ws.Range("E19").Value = "DigiKey"
$ws.Hyperlinks.Add($ws.Range("E19"), "https://www.digikey.com/en/products/detail/molex/0510210100/26427") | Out-Null
$ws.Range("E19").HorizontalAlignment = -4108

$ws.Range("I19").Value = "Fall 2024"

# ---------------------------------------------------------------------
# Row 20 - Raspberry pi
# ---------------------------------------------------------------------
$ws.Range("A20").Value = "Raspberry pi "
$ws.Range("A20").HorizontalAlignment = -4108

$ws.Range("B20").Value = 2
$ws.Range("B20").HorizontalAlignment = -4108

$ws.Range("C20").Value = 35
$ws.Range("C20").NumberFormat = $CURRENCY0
$ws.Range("C20").HorizontalAlignment = -4108

$ws.Range("D20").Value = 70
$ws.Range("D20").NumberFormat = $CURRENCY0
$ws.Range("D20").HorizontalAlignment = -4108

$ws.Range("E20").Value = "Digikey"
$ws.Hyperlinks.Add($ws.Range("E20"), "https://www.digikey.com/en/products/detail/raspberry-pi/SC0194/13624703") | Out-Null
$ws.Range("E20").HorizontalAlignment = -4108

$ws.Range("I20").Value = "Fall 2024"

# ---------------------------------------------------------------------
# Clear the scratch seeding cells - they must not remain in the sheet.
# ---------------------------------------------------------------------
$ws.Range($ws.Cells.Item($seedRow, 1), $ws.Cells.Item($seedRow, $seedStrings.Count)).ClearContents()

# ---------------------------------------------------------------------
# View: scroll window + selection, matching the saved sheet view state
# ---------------------------------------------------------------------
$ws.Range("G18").Select()
$excel.ActiveWindow.ScrollRow = 5

Write-Output "done"
